$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 162, shifting existing rows 162-180 down to 163-181
$ws.Range("A162:R162").EntireRow.Insert(-4121)

# Populate the newly inserted row 162 with the new record
$ws.Range("A162").Value = 10
$ws.Range("B162").Value = "Vega Modelo de Temuco"
$ws.Range("C162").Value = "La Araucanía"
$ws.Range("D162").Value = 44776
$ws.Range("E162").Value = 9
$ws.Range("F162").Value = 100114007
$ws.Range("G162").Value = "Jengibre"
$ws.Range("H162").Value = "Sin especificar"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 30
$ws.Range("K162").Value = 20000
$ws.Range("L162").Value = 20000
$ws.Range("M162").Value = 20000
$ws.Range("N162").Value = "`$/caja 13 kilos"
$ws.Range("O162").Value = "Perú"
$ws.Range("P162").Value = 1538
$ws.Range("Q162").Value = 13
$ws.Range("R162").Value = "Hortaliza"
